$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (pushes existing rows 2-6 down to 3-7)
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "2021-01"
$ws.Range("B2").Value = "'233"

# Insert a new row at row 4 (before the "2022-02" row, now at row 4)
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "2022-01"
$ws.Range("B4").Value = "'789"

# Update the "2023-01" row's consumption value (now at row 6)
$ws.Range("B6").Value = "'123123123"

# Append two new rows at the end
$ws.Range("A9").Value = "2023-10"
$ws.Range("B9").Value = "'550"

$ws.Range("A10").Value = "2023-11"
$ws.Range("B10").Value = "'500"
